$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Payment Screenshot" header column (S1), matching the style of the other headers ---
$ws.Range("S1").Value = "Payment Screenshot"
$ws.Range("R1").Copy() | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null

# --- Add blank text placeholder cells in column S for the existing data rows (S2:S7) ---
$ws.Range("S2:S7").Value = "'"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("S2:S7").PasteSpecial(-4122) | Out-Null

# --- Row 7: Pincode (I7) switches from text to a genuine number ---
$ws.Range("I7").Value = 121004

# --- Row 8: new order row. Seed it from row 7 (same blank/text/number cell pattern), then
#     overwrite the cells that actually differ for this new order. ---
$ws.Range("A7:R7").Copy() | Out-Null
$ws.Range("A8:R8").Select() | Out-Null
$ws.Paste() | Out-Null

$ws.Range("H8").Value = "jeevan"
$ws.Range("L8").Value = "ORD20250731165620"
$ws.Range("N8").Value = "Light Green Dress (₹250 x 1), Mala Set of 8 (₹240 x 1)"
$ws.Range("O8").Value = 549
$ws.Range("P8").Value = "2025-07-31 16:56"
$ws.Range("Q8").Value = 490
$ws.Range("R8").Value = 59

# A8:F8 and J8:K8 need to stay as blank text cells (same pattern as the other rows)
$ws.Range("A8:F8").Value = "'"
$ws.Range("J8:K8").Value = "'"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A8:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("J8:K8").PasteSpecial(-4122) | Out-Null

# I8 keeps its pincode text ('121004') without the leftover quote-prefix style
$ws.Range("I8").Value = "'121004"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null

# S8: payment screenshot filename for the new order
$ws.Range("S8").Value = "static/payments/ORD20250731165620.jpeg"

